$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Reword the two "Y:.../N:..." notes to put the N: part on its own line
# (Alt+Enter line break inside the cell instead of "; "/",").
$ws.Range("G14").Value = "Y:是" + [char]10 + "N:否"
$ws.Range("G14").WrapText = $true
$ws.Range("G16").Value = "Y:同意" + [char]10 + "N:不同意"

# Rows auto-grow to fit the now two-line wrapped text.
$ws.Rows.Item(14).RowHeight = 32.4
$ws.Rows.Item(16).RowHeight = 32.4

# Move the view: no more frozen-ish scroll to A4, and the remembered
# selection is now G17 instead of G29.
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("G17").Select()
